# Update "取得日時" (retrieved-at) timestamps on the "ランサーズ" sheet.
# All existing data rows (2-8) get their timestamp bumped to the new
# scrape time: 2026-01-19 02:06:23 (JST), per commit "Append: 2026-01-19 02:06 JST".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2026-01-19 02:06:23"

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 2) {
    $lastRow = 8
}

for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 1).Value = $newTimestamp
}
